$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 2 (pushes existing rows 2-5 down to 4-7)
$ws.Rows("2:3").Insert()
# The inserted rows picked up formatting from the row above (header); clear it
# so the new data rows match the unstyled look of the other data rows.
$ws.Rows("2:3").ClearFormats()

# ---- New row 2: FIFA World Cup Qualifiers - Asia (Oman vs Qatar) ----
$ws.Range("A2").Value2 = "FIFA World Cup Qualifiers - Asia"
# Dates like "2025-10-08" get auto-parsed into date serials by Value2;
# force text storage, write, then clear the temporary number format.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value2 = "2025-10-08"
$ws.Range("B2").ClearFormats()
$ws.Range("C2").Value2 = "12:00:00"
$ws.Range("D2").Value2 = "Oman"
$ws.Range("E2").Value2 = "Qatar"
$ws.Range("F2").Value2 = 6.4
$ws.Range("G2").Value2 = 7.8
$ws.Range("H2").Value2 = 1.65
$ws.Range("I2").Value2 = 1.73
$ws.Range("J2").Value2 = 3.7
$ws.Range("K2").Value2 = 3.9
$ws.Range("L2").Value2 = 1.01
$ws.Range("M2").Value2 = 1.08
$ws.Range("N2").Value2 = 2.82
$ws.Range("O2").Value2 = 1.44
$ws.Range("P2").Value2 = 1.62
$ws.Range("Q2").Value2 = 2.32
$ws.Range("R2").Value2 = 1.22
$ws.Range("S2").Value2 = 4.5
$ws.Range("T2").Value2 = 2.22
$ws.Range("U2").Value2 = 1.68
$ws.Range("V2").Value2 = 2.12
$ws.Range("W2").Value2 = 1.12
$ws.Range("X2").Value2 = 11
$ws.Range("Y2").Value2 = 6.6
$ws.Range("Z2").Value2 = 9.199999999999999
$ws.Range("AA2").Value2 = 18
$ws.Range("AB2").Value2 = 18.5
$ws.Range("AC2").Value2 = 9
$ws.Range("AD2").Value2 = 11
$ws.Range("AE2").Value2 = 23
$ws.Range("AF2").Value2 = 60
$ws.Range("AG2").Value2 = 1000
$ws.Range("AH2").Value2 = 1000
$ws.Range("AI2").Value2 = 60
$ws.Range("AJ2").Value2 = 320
$ws.Range("AK2").Value2 = 180
$ws.Range("AL2").Value2 = 180
$ws.Range("AM2").Value2 = 280
$ws.Range("AN2").Value2 = 1000
$ws.Range("AO2").Value2 = 16

# ---- New row 3: FIFA World Cup Qualifiers - Asia (Indonesia vs Saudi Arabia) ----
$ws.Range("A3").Value2 = "FIFA World Cup Qualifiers - Asia"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value2 = "2025-10-08"
$ws.Range("B3").ClearFormats()
$ws.Range("C3").Value2 = "14:15:00"
$ws.Range("D3").Value2 = "Indonesia"
$ws.Range("E3").Value2 = "Saudi Arabia"
$ws.Range("F3").Value2 = 7.8
$ws.Range("G3").Value2 = 10
$ws.Range("H3").Value2 = 1.49
$ws.Range("I3").Value2 = 1.55
$ws.Range("J3").Value2 = 4.2
$ws.Range("K3").Value2 = 4.5
$ws.Range("L3").Value2 = 1.42
$ws.Range("M3").Value2 = 1.07
$ws.Range("N3").Value2 = 3.15
$ws.Range("O3").Value2 = 1.38
$ws.Range("P3").Value2 = 1.74
$ws.Range("Q3").Value2 = 2.12
$ws.Range("R3").Value2 = 1.27
$ws.Range("S3").Value2 = 3.95
$ws.Range("T3").Value2 = 2.06
$ws.Range("U3").Value2 = 1.55
$ws.Range("V3").Value2 = 2.52
$ws.Range("W3").Value2 = 1.09
$ws.Range("X3").Value2 = 15.5
$ws.Range("Y3").Value2 = 6.6
$ws.Range("Z3").Value2 = 8
$ws.Range("AA3").Value2 = 13.5
$ws.Range("AB3").Value2 = 24
$ws.Range("AC3").Value2 = 10.5
$ws.Range("AD3").Value2 = 11
$ws.Range("AE3").Value2 = 19.5
$ws.Range("AF3").Value2 = 100
$ws.Range("AG3").Value2 = 36
$ws.Range("AH3").Value2 = 34
$ws.Range("AI3").Value2 = 55
$ws.Range("AJ3").Value2 = 1000
$ws.Range("AK3").Value2 = 230
$ws.Range("AL3").Value2 = 210
$ws.Range("AM3").Value2 = 1000
$ws.Range("AN3").Value2 = 1000
$ws.Range("AO3").Value2 = 10.5

# ---- Row 4 (formerly row 2): Colombian Primera A - updated odds ----
$ws.Range("G4").Value2 = 4.4
$ws.Range("H4").Value2 = 2.22
$ws.Range("I4").Value2 = 2.54
$ws.Range("J4").Value2 = 2.82
$ws.Range("K4").Value2 = 3.7
$ws.Range("P4").Value2 = 1.73
$ws.Range("Q4").Value2 = 2.1

# ---- Row 5 (formerly row 3): Brazilian Serie A - updated odds ----
$ws.Range("F5").Value2 = 1.66
$ws.Range("G5").Value2 = 1.7
$ws.Range("H5").Value2 = 7.2
$ws.Range("J5").Value2 = 3.7
$ws.Range("K5").Value2 = 3.9
$ws.Range("P5").Value2 = 1.71
